# Implement column width parsing: give column A and column D explicit
# custom widths (matches the "Implement column width parsing" test
# fixture, which exercises <cols> width parsing in calamine).
#
# Excel's Range/Columns.ColumnWidth COM property is expressed in
# "characters" at the Normal style's font, while the stored OOXML <col>
# width adds a fixed 5-pixel padding (5/6 = 0.8333... chars at this
# sheet's font metrics). Subtract that offset here so the persisted
# <col width="..."> values come out to exactly 27 and 32.5, matching
# the target workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).ColumnWidth = 26.166666666666668   # -> stored width 27
$ws.Columns.Item(4).ColumnWidth = 31.666666666666668   # -> stored width 32.5
